$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.979.79"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.17%  "
$ws.Range("D3").Value = "'2.734.22"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.89%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'569.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.59%  "
$ws.Range("E6").Value = "  -0.32%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "'0.596"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.14%  "
$ws.Range("E9").Value = "  -2.33%  "
$ws.Range("E10").Value = "  +4.37%  "
$ws.Range("D11").Value = "'5.68"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.37%  "
$ws.Range("E12").Value = "  -2.03%  "
$ws.Range("D13").Value = "'3.216.72"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.94%  "
$ws.Range("D14").Value = "'26.65"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.11%  "
$ws.Range("D15").Value = "'63.583.93"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.66%  "
$ws.Range("E16").Value = "  -2.58%  "
$ws.Range("D17").Value = "'2.735.46"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.99%  "
$ws.Range("D18").Value = "'12.05"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.70%  "
$ws.Range("E19").Value = "  -2.45%  "
$ws.Range("D20").Value = "'353.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.08%  "
$ws.Range("D21").Value = "'6.55"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.88%  "
$ws.Range("D22").Value = "'0.996"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.35%  "
$ws.Range("D23").Value = "'0.521"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.81%  "
$ws.Range("D24").Value = "'64.28"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.80%  "
$ws.Range("E25").Value = "  -1.00%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("D27").Value = "'8.41"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.63%  "
$ws.Range("D28").Value = "'0.0₃0909"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.35%  "
$ws.Range("D29").Value = "'1.97"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.08%  "
$ws.Range("B30").Value = "Aptos"
$ws.Range("C30").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D30").Value = "'7.21"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.22%  "
$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").Value = "'1.35"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.36%  "
$ws.Range("D32").Value = "'162.94"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.68%  "
$ws.Range("E33").Value = "  -2.12%  "
$ws.Range("E34").Value = "  -1.89%  "
$ws.Range("E35").Value = "  +0.17%  "
$ws.Range("E37").Value = "  -0.74%  "
$ws.Range("D38").Value = "'0.987"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.41%  "
$ws.Range("D39").Value = "'348.72"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.48%  "
$ws.Range("D40").Value = "'6.28"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.37%  "
$ws.Range("D41").Value = "'4.12"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.70%  "
$ws.Range("E42").Value = "  -1.69%  "
$ws.Range("E43").Value = "  +0.29%  "
$ws.Range("E44").Value = "  -3.30%  "
$ws.Range("E45").Value = "  -2.70%  "
$ws.Range("D46").Value = "'134.59"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.56%  "
$ws.Range("D47").Value = "'0.622"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.62%  "
$ws.Range("E48").Value = "  -1.74%  "
$ws.Range("D49").Value = "'0.0248"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.05%  "
$ws.Range("D50").Value = "'0.998"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.09%  "
$ws.Range("D51").Value = "'11.02"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.26%  "
